$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list on Sat Nov 18 21:45:41 UTC 2023 with GitHub Actions
# Refresh price (col D) and 1h volume-change (col E) snapshot values.
# A few D-column values look like plain decimals (e.g. "1.00", "22.22"); force
# NumberFormat to Text first so Excel keeps them as literal strings instead of
# auto-coercing to numbers (which would also drop significant trailing zeros).

$ws.Range("D2").Value = "36.598.35"
$ws.Range("E2").Value = "  +0.56%  "
$ws.Range("D3").Value = "1.960.85"
$ws.Range("E3").Value = "  +0.88%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.47"
$ws.Range("E5").Value = "  +0.62%  "
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.83"
$ws.Range("E7").Value = "  +0.71%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.376"
$ws.Range("E9").Value = "  +3.27%  "
$ws.Range("E10").Value = "  -2.28%  "
$ws.Range("E11").Value = "  -0.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.22"
$ws.Range("E12").Value = "  +3.25%  "
$ws.Range("D13").Value = "2.248.65"
$ws.Range("E13").Value = "  +0.90%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.826"
$ws.Range("E14").Value = "  +0.73%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "13.70"
$ws.Range("E15").Value = "  +0.80%  "
$ws.Range("E16").Value = "  +0.67%  "
$ws.Range("D17").Value = "1.963.23"
$ws.Range("E17").Value = "  -0.88%  "
$ws.Range("D18").Value = "36.489.79"
$ws.Range("E18").Value = "  +0.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.87"
$ws.Range("E19").Value = "  +0.54%  "
$ws.Range("D20").Value = "0.0₃0857"
$ws.Range("E20").Value = "  -0.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "228.30"
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.05"
$ws.Range("E22").Value = "  +0.39%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("E24").Value = "  +0.99%  "
$ws.Range("E25").Value = "  +3.30%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.140"
$ws.Range("E26").Value = "  +7.96%  "
$ws.Range("E27").Value = "  +0.15%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "160.07"
$ws.Range("E28").Value = "  -0.92%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.41"
$ws.Range("E29").Value = "  -0.28%  "
$ws.Range("E31").Value = "  +1.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.73"
$ws.Range("E32").Value = "  +1.55%  "
$ws.Range("E33").Value = "  -1.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.28"
$ws.Range("E34").Value = "  +0.56%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("E36").Value = "  +5.93%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.36"
$ws.Range("E37").Value = "  +11.03%  "
$ws.Range("B38").Value = "THORChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.96"
$ws.Range("E38").Value = "  -4.83%  "
$ws.Range("E39").Value = "  -0.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0984"
$ws.Range("E40").Value = "  +0.56%  "
$ws.Range("E41").Value = "  +1.08%  "
$ws.Range("E42").Value = "  +0.39%  "
$ws.Range("E43").Value = "  +0.88%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.97"
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("D45").Value = "1.365.29"
$ws.Range("E45").Value = "  +0.94%  "
$ws.Range("E46").Value = "  +0.62%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "87.70"
$ws.Range("E47").Value = "  -0.02%  "
$ws.Range("E48").Value = "  +0.56%  "
$ws.Range("E49").Value = "  +0.58%  "
$ws.Range("D50").Value = "2.140.19"
$ws.Range("E50").Value = "  +1.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "43.64"
$ws.Range("E51").Value = "  -3.55%  "
